# Generate Report for Handoff
# Adds two new localized files (75e92f29-... and a98889f4-...) as new rows
# to the "Overview", "zh-cn" and "de-de" worksheets of the localization
# status report, and resizes the corresponding tables/ranges to include them.

$wb = $excel.ActiveWorkbook

$file1Guid = "75e92f29-baf8-400f-bf29-5730c23c8f4b"
$file2Guid = "a98889f4-18d6-4bd6-a462-69f1d0248861"

$file1Name = "$file1Guid.md"
$file2Name = "$file2Guid.md"

$file1Path = "e2e\$file1Guid.md"
$file2Path = "e2e\$file2Guid.md"

$file1ZhXlf = "$file1Guid.1a1612434d1058199b93f8f85dca8843de3d3f6f.zh-cn.xlf"
$file2ZhXlf = "$file2Guid.a8c7c336b80e01e7fb6de4bda5e29a4a7c1b2614.zh-cn.xlf"
$file1DeXlf = "$file1Guid.1a1612434d1058199b93f8f85dca8843de3d3f6f.de-de.xlf"
$file2DeXlf = "$file2Guid.a8c7c336b80e01e7fb6de4bda5e29a4a7c1b2614.de-de.xlf"

$overviewDate = "2016-08-17 22:39:45"
$zhHandoffDate = "2016-08-17 22:39:40"
$deHandoffDate = "2016-08-17 22:39:45"

$file1ZhUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/1a1612434d1058199b93f8f85dca8843de3d3f6f/$file1Path"
$file2ZhUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a8c7c336b80e01e7fb6de4bda5e29a4a7c1b2614/$file2Path"
$file1DeUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/1a1612434d1058199b93f8f85dca8843de3d3f6f/$file1Path"
$file2DeUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a8c7c336b80e01e7fb6de4bda5e29a4a7c1b2614/$file2Path"
$file1OvUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/1a1612434d1058199b93f8f85dca8843de3d3f6f/$file1Path"
$file2OvUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a8c7c336b80e01e7fb6de4bda5e29a4a7c1b2614/$file2Path"

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Cells.Item(4,1).Value2 = $file1Name
$wsOverview.Cells.Item(4,2).Value2 = $file1Path
$wsOverview.Cells.Item(4,3).Value2 = ".md"
$wsOverview.Cells.Item(4,5).Value2 = "Ready for handoff"
$wsOverview.Cells.Item(4,6).Value2 = "Ready for handoff"
$wsOverview.Cells.Item(4,7).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsOverview.Cells.Item(4,7).Value2 = $overviewDate

$wsOverview.Cells.Item(5,1).Value2 = $file2Name
$wsOverview.Cells.Item(5,2).Value2 = $file2Path
$wsOverview.Cells.Item(5,3).Value2 = ".md"
$wsOverview.Cells.Item(5,5).Value2 = "Ready for handoff"
$wsOverview.Cells.Item(5,6).Value2 = "Ready for handoff"
$wsOverview.Cells.Item(5,7).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsOverview.Cells.Item(5,7).Value2 = $overviewDate

$wsOverview.Hyperlinks.Add($wsOverview.Range("B4"), $file1OvUrl, "", "", $file1Path) | Out-Null
$wsOverview.Hyperlinks.Add($wsOverview.Range("B5"), $file2OvUrl, "", "", $file2Path) | Out-Null

$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.Resize($wsOverview.Range("A1:G5"))

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Cells.Item(4,1).Value2 = $file1Name
$wsZh.Cells.Item(4,2).Value2 = ".md"
$wsZh.Cells.Item(4,3).Value2 = "Ready for handoff"
$wsZh.Cells.Item(4,4).Value2 = "e2e"
$wsZh.Cells.Item(4,5).Value2 = "ht"
$wsZh.Cells.Item(4,6).Value2 = "'False"
$wsZh.Cells.Item(4,7).Value2 = $file1ZhXlf
$wsZh.Cells.Item(4,8).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Cells.Item(4,8).Value2 = $zhHandoffDate
$wsZh.Cells.Item(4,11).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Cells.Item(4,11).Value2 = "0001-01-01 00:00:00"
$wsZh.Cells.Item(4,13).Value2 = "'True"
$wsZh.Cells.Item(4,15).Value2 = "'False"

$wsZh.Cells.Item(5,1).Value2 = $file2Name
$wsZh.Cells.Item(5,2).Value2 = ".md"
$wsZh.Cells.Item(5,3).Value2 = "Ready for handoff"
$wsZh.Cells.Item(5,4).Value2 = "e2e"
$wsZh.Cells.Item(5,5).Value2 = "ht"
$wsZh.Cells.Item(5,6).Value2 = "'False"
$wsZh.Cells.Item(5,7).Value2 = $file2ZhXlf
$wsZh.Cells.Item(5,8).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Cells.Item(5,8).Value2 = $zhHandoffDate
$wsZh.Cells.Item(5,11).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Cells.Item(5,11).Value2 = "0001-01-01 00:00:00"
$wsZh.Cells.Item(5,13).Value2 = "'True"
$wsZh.Cells.Item(5,15).Value2 = "'False"

$wsZh.Hyperlinks.Add($wsZh.Range("A4"), $file1ZhUrl, "", "", $file1Name) | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("A5"), $file2ZhUrl, "", "", $file2Name) | Out-Null

$loZh = $wsZh.ListObjects.Item(1)
$loZh.Resize($wsZh.Range("A1:P5"))

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Cells.Item(4,1).Value2 = $file1Name
$wsDe.Cells.Item(4,2).Value2 = ".md"
$wsDe.Cells.Item(4,3).Value2 = "Ready for handoff"
$wsDe.Cells.Item(4,4).Value2 = "e2e"
$wsDe.Cells.Item(4,5).Value2 = "ht"
$wsDe.Cells.Item(4,6).Value2 = "'False"
$wsDe.Cells.Item(4,7).Value2 = $file1DeXlf
$wsDe.Cells.Item(4,8).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Cells.Item(4,8).Value2 = $deHandoffDate
$wsDe.Cells.Item(4,11).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Cells.Item(4,11).Value2 = "0001-01-01 00:00:00"
$wsDe.Cells.Item(4,13).Value2 = "'True"
$wsDe.Cells.Item(4,15).Value2 = "'False"

$wsDe.Cells.Item(5,1).Value2 = $file2Name
$wsDe.Cells.Item(5,2).Value2 = ".md"
$wsDe.Cells.Item(5,3).Value2 = "Ready for handoff"
$wsDe.Cells.Item(5,4).Value2 = "e2e"
$wsDe.Cells.Item(5,5).Value2 = "ht"
$wsDe.Cells.Item(5,6).Value2 = "'False"
$wsDe.Cells.Item(5,7).Value2 = $file2DeXlf
$wsDe.Cells.Item(5,8).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Cells.Item(5,8).Value2 = $deHandoffDate
$wsDe.Cells.Item(5,11).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Cells.Item(5,11).Value2 = "0001-01-01 00:00:00"
$wsDe.Cells.Item(5,13).Value2 = "'True"
$wsDe.Cells.Item(5,15).Value2 = "'False"

$wsDe.Hyperlinks.Add($wsDe.Range("A4"), $file1DeUrl, "", "", $file1Name) | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("A5"), $file2DeUrl, "", "", $file2Name) | Out-Null

$loDe = $wsDe.ListObjects.Item(1)
$loDe.Resize($wsDe.Range("A1:P5"))
